# Variables.xlsx: "Added 1 Variable" — insert a new variable row "HSD010"
# before the WHD110/WHD120/BMXWT block (old row 75 -> new row 76), and
# move the active selection to F66 (matching the new sheetView selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 75, pushing WHD110/WHD120/BMXWT down by one.
$ws.Rows("75:75").Insert()

# New variable name goes into the freshly inserted row.
$ws.Range("A75").Value = "HSD010"

# Reflect the author's final cursor position in the saved view.
[void]$ws.Range("F66").Select()
